$d = $word.ActiveDocument
$lines = @(
    'Head west on SW 1st St toward SW 38th Ave',
    'Turn right onto SW 38th Ave',
    'Turn left onto W Broward Blvd',
    'Turn right onto US-441 NN State Rd 7',
    'Exit on the left onto FL-838 WW Sunrise Blvd',
    'Turn right to merge onto Florida''s Turnpike toward Orlando',
    'Merge onto I-75 N',
    'Take exit 435 to merge onto I-10 W toward Tallahassee',
    'Take the I-12 W exit toward HammondBaton Rouge',
    'Continue onto I-12 W',
    'Merge onto I-10 W',
    'Keep left to stay on I-10 W',
    'Keep left to stay on I-10 W',
    'Keep left to stay on I-10 W',
    'Sharp left to stay on I-10 W',
    'Keep left at the fork to stay on I-10 W',
    'Exit onto Interstate 10 Access Rd',
    'Keep right to stay on Interstate 10 Access Rd',
    'Continue straight to stay on Interstate 10 Access Rd',
    'Continue straight to stay on Interstate 10 Access Rd',
    'Continue straight',
    'Merge onto I-10 W',
    'Keep left to stay on I-10 W',
    'Keep right to stay on I-10 W',
    'Keep right to stay on I-10 W',
    'Keep right to stay on I-10 W',
    'Continue onto I-10 W',
    'Keep left to stay on I-10 W',
    'Keep right',
    'Keep right',
    'Slight left onto I-10 Frontage Rd',
    'Continue onto US-90 W',
    'Keep left',
    'Continue straight',
    'Take the ramp on the left onto I-10 W',
    'Take exit 587 to merge onto I-10 Frontage Rd',
    'Turn right onto E Charles William Anderson LoopN Loop 1604 E',
    'Keep right to continue on TX-1604 LoopN Loop 1604 E',
    'Take the exit toward I-10 WUS-87 N',
    'Keep right at the fork, follow signs for I-10 WUS-87 NEl Paso and merge onto I-10 WUS-87 N',
    'Keep left at the fork to stay on I-10 W',
    'Take exit 199 for I-8 W toward San Diego',
    'Continue onto I-8 W',
    'Take exit 14B to merge onto CA-125 S',
    'Keep left at the fork, follow signs for CA-94 W',
    'Continue onto CA-94 W',
    'Take exit 2 on the left for CA-15 SEscondido Fwy',
    'Merge onto CA-15I-15 SEscondido Fwy',
    'Keep right to continue on I-15 S, follow signs for I-5 N',
    'Take exit 1B to merge onto I-5 N',
    'Take exit 14A for CA-75 toward Coronado',
    'Continue onto CA-75 SSan Diego – Coronado Bridge',
    'Turn left onto Orange Ave',
    'Turn left at the 1st cross street onto 4th St',
    'Turn right onto Pomona Ave',
    'At the traffic circle, take the 3rd exit and stay on Pomona Ave',
    'Turn right onto Glorietta Blvd',
    'Keep left to continue on Pomona Ave',
    'Turn left onto Strand Way'
)
$joined = [string]::Join([char]13, $lines)
$rng = $d.Content
$rng.Collapse(0)
$rng.InsertAfter([char]13 + $joined)
Write-Host "Appended $($lines.Count) direction paragraphs"
